$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bill Russell per-game averages: an updated Kaggle CSV export reshuffled
# which season stat line lives on which worksheet row. Column B
# (season_ending_year) must stay text-typed like the original cells.

# Force column B to keep its original text type for every row whose
# season changes (rows 6 and 14 are untouched by this edit).
$bRows = @(2,3,4,5,7,8,9,10,11,12,13)
foreach ($r in $bRows) {
    $ws.Cells.Item($r, 2).NumberFormat = "@"
}

# Row 2
$ws.Range("A2").Value = 2735
$ws.Range("B2").Value = '1965'
$ws.Range("G2").Value = 31
$ws.Range("H2").Value = 9
$ws.Range("K2").Value = 78
$ws.Range("M2").Value = 44.4
$ws.Range("N2").Value = 5.5
$ws.Range("O2").Value = 12.6
$ws.Range("P2").Value = 0.438
$ws.Range("T2").Value = 5.5
$ws.Range("U2").Value = 12.6
$ws.Range("V2").Value = 0.438
$ws.Range("W2").Value = 0.438
$ws.Range("X2").Value = 3.1
$ws.Range("Y2").Value = 5.5
$ws.Range("Z2").Value = 0.573
$ws.Range("AC2").Value = 24.1
$ws.Range("AD2").Value = 5.3
$ws.Range("AH2").Value = 2.6
$ws.Range("AI2").Value = 14.1
$ws.Range("AJ2").Value = '1964-65'
$ws.Range("AM2").Value = 1965

# Row 3
$ws.Range("A3").Value = 1760
$ws.Range("B3").Value = '1957'
$ws.Range("G3").Value = 23
$ws.Range("H3").Value = 1
$ws.Range("K3").Value = 48
$ws.Range("M3").Value = 35.3
$ws.Range("N3").Value = 5.8
$ws.Range("O3").Value = 13.5
$ws.Range("P3").Value = 0.427
$ws.Range("T3").Value = 5.8
$ws.Range("U3").Value = 13.5
$ws.Range("V3").Value = 0.427
$ws.Range("W3").Value = 0.427
$ws.Range("X3").Value = 3.2
$ws.Range("Y3").Value = 6.4
$ws.Range("Z3").Value = 0.492
$ws.Range("AC3").Value = 19.6
$ws.Range("AD3").Value = 1.8
$ws.Range("AH3").Value = 3
$ws.Range("AI3").Value = 14.7
$ws.Range("AJ3").Value = '1956-57'
$ws.Range("AM3").Value = 1957

# Row 4
$ws.Range("A4").Value = 2219
$ws.Range("B4").Value = '1961'
$ws.Range("G4").Value = 27
$ws.Range("H4").Value = 5
$ws.Range("M4").Value = 44.3
$ws.Range("N4").Value = 6.8
$ws.Range("O4").Value = 16
$ws.Range("P4").Value = 0.426
$ws.Range("T4").Value = 6.8
$ws.Range("U4").Value = 16
$ws.Range("V4").Value = 0.426
$ws.Range("W4").Value = 0.426
$ws.Range("X4").Value = 3.3
$ws.Range("Y4").Value = 6
$ws.Range("Z4").Value = 0.55
$ws.Range("AC4").Value = 23.9
$ws.Range("AD4").Value = 3.4
$ws.Range("AH4").Value = 2
$ws.Range("AI4").Value = 16.9
$ws.Range("AJ4").Value = '1960-61'
$ws.Range("AM4").Value = 1961

# Row 5
$ws.Range("A5").Value = 2097
$ws.Range("B5").Value = '1960'
$ws.Range("G5").Value = 26
$ws.Range("H5").Value = 4
$ws.Range("K5").Value = 74
$ws.Range("M5").Value = 42.5
$ws.Range("N5").Value = 7.5
$ws.Range("O5").Value = 16.1
$ws.Range("P5").Value = 0.467
$ws.Range("T5").Value = 7.5
$ws.Range("U5").Value = 16.1
$ws.Range("V5").Value = 0.467
$ws.Range("W5").Value = 0.467
$ws.Range("X5").Value = 3.2
$ws.Range("Y5").Value = 5.3
$ws.Range("Z5").Value = 0.612
$ws.Range("AC5").Value = 24
$ws.Range("AD5").Value = 3.7
$ws.Range("AH5").Value = 2.8
$ws.Range("AI5").Value = 18.2
$ws.Range("AJ5").Value = '1959-60'
$ws.Range("AM5").Value = 1960

# Row 7
$ws.Range("A7").Value = 2865
$ws.Range("B7").Value = '1966'
$ws.Range("G7").Value = 32
$ws.Range("H7").Value = 10
$ws.Range("K7").Value = 78
$ws.Range("M7").Value = 43.4
$ws.Range("N7").Value = 5
$ws.Range("O7").Value = 12.1
$ws.Range("P7").Value = 0.415
$ws.Range("T7").Value = 5
$ws.Range("U7").Value = 12.1
$ws.Range("V7").Value = 0.415
$ws.Range("W7").Value = 0.415
$ws.Range("X7").Value = 2.9
$ws.Range("Y7").Value = 5.2
$ws.Range("Z7").Value = 0.551
$ws.Range("AC7").Value = 22.8
$ws.Range("AD7").Value = 4.8
$ws.Range("AH7").Value = 2.8
$ws.Range("AI7").Value = 12.9
$ws.Range("AJ7").Value = '1965-66'
$ws.Range("AM7").Value = 1966

# Row 8
$ws.Range("A8").Value = 2332
$ws.Range("B8").Value = '1962'
$ws.Range("G8").Value = 28
$ws.Range("H8").Value = 6
$ws.Range("K8").Value = 76
$ws.Range("M8").Value = 45.2
$ws.Range("N8").Value = 7.6
$ws.Range("O8").Value = 16.6
$ws.Range("P8").Value = 0.457
$ws.Range("T8").Value = 7.6
$ws.Range("U8").Value = 16.6
$ws.Range("V8").Value = 0.457
$ws.Range("W8").Value = 0.457
$ws.Range("X8").Value = 3.8
$ws.Range("Y8").Value = 6.3
$ws.Range("Z8").Value = 0.595
$ws.Range("AC8").Value = 23.6
$ws.Range("AD8").Value = 4.5
$ws.Range("AH8").Value = 2.7
$ws.Range("AI8").Value = 18.9
$ws.Range("AJ8").Value = '1960-62'
$ws.Range("AM8").Value = 1962

# Row 9
$ws.Range("A9").Value = 3513
$ws.Range("B9").Value = '1969'
$ws.Range("G9").Value = 35
$ws.Range("H9").Value = 13
$ws.Range("K9").Value = 77
$ws.Range("M9").Value = 42.7
$ws.Range("N9").Value = 3.6
$ws.Range("O9").Value = 8.4
$ws.Range("P9").Value = 0.433
$ws.Range("T9").Value = 3.6
$ws.Range("U9").Value = 8.4
$ws.Range("V9").Value = 0.433
$ws.Range("W9").Value = 0.433
$ws.Range("X9").Value = 2.6
$ws.Range("Y9").Value = 5
$ws.Range("Z9").Value = 0.526
$ws.Range("AC9").Value = 19.3
$ws.Range("AD9").Value = 4.9
$ws.Range("AH9").Value = 3
$ws.Range("AI9").Value = 9.9
$ws.Range("AJ9").Value = '1968-69'
$ws.Range("AK9").Value = 'No'
$ws.Range("AM9").Value = 1969

# Row 10
$ws.Range("A10").Value = 2600
$ws.Range("B10").Value = '1964'
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = 8
$ws.Range("M10").Value = 44.6
$ws.Range("N10").Value = 6
$ws.Range("O10").Value = 13.8
$ws.Range("P10").Value = 0.433
$ws.Range("T10").Value = 6
$ws.Range("U10").Value = 13.8
$ws.Range("V10").Value = 0.433
$ws.Range("W10").Value = 0.433
$ws.Range("X10").Value = 3
$ws.Range("Y10").Value = 5.5
$ws.Range("Z10").Value = 0.55
$ws.Range("AC10").Value = 24.7
$ws.Range("AD10").Value = 4.7
$ws.Range("AH10").Value = 2.4
$ws.Range("AI10").Value = 15
$ws.Range("AJ10").Value = '1963-64'
$ws.Range("AM10").Value = 1964

# Row 11
$ws.Range("A11").Value = 2467
$ws.Range("B11").Value = '1963'
$ws.Range("G11").Value = 29
$ws.Range("H11").Value = 7
$ws.Range("K11").Value = 78
$ws.Range("M11").Value = 44.9
$ws.Range("N11").Value = 6.6
$ws.Range("O11").Value = 15.2
$ws.Range("P11").Value = 0.432
$ws.Range("T11").Value = 6.6
$ws.Range("U11").Value = 15.2
$ws.Range("V11").Value = 0.432
$ws.Range("W11").Value = 0.432
$ws.Range("X11").Value = 3.7
$ws.Range("Y11").Value = 6.6
$ws.Range("Z11").Value = 0.555
$ws.Range("AC11").Value = 23.6
$ws.Range("AD11").Value = 4.5
$ws.Range("AH11").Value = 2.4
$ws.Range("AI11").Value = 16.8
$ws.Range("AJ11").Value = '1962-63'
$ws.Range("AM11").Value = 1963

# Row 12
$ws.Range("A12").Value = 1875
$ws.Range("B12").Value = '1958'
$ws.Range("G12").Value = 24
$ws.Range("H12").Value = 2
$ws.Range("K12").Value = 69
$ws.Range("M12").Value = 38.3
$ws.Range("N12").Value = 6.6
$ws.Range("O12").Value = 15
$ws.Range("P12").Value = 0.442
$ws.Range("T12").Value = 6.6
$ws.Range("U12").Value = 15
$ws.Range("V12").Value = 0.442
$ws.Range("W12").Value = 0.442
$ws.Range("X12").Value = 3.3
$ws.Range("Y12").Value = 6.4
$ws.Range("Z12").Value = 0.519
$ws.Range("AC12").Value = 22.7
$ws.Range("AD12").Value = 2.9
$ws.Range("AH12").Value = 2.6
$ws.Range("AI12").Value = 16.6
$ws.Range("AJ12").Value = '1957-58'
$ws.Range("AM12").Value = 1958

# Row 13
$ws.Range("A13").Value = 3150
$ws.Range("B13").Value = '1968'
$ws.Range("G13").Value = 34
$ws.Range("H13").Value = 12
$ws.Range("K13").Value = 78
$ws.Range("M13").Value = 37.9
$ws.Range("N13").Value = 4.7
$ws.Range("O13").Value = 11
$ws.Range("P13").Value = 0.425
$ws.Range("T13").Value = 4.7
$ws.Range("U13").Value = 11
$ws.Range("V13").Value = 0.425
$ws.Range("W13").Value = 0.425
$ws.Range("X13").Value = 3.2
$ws.Range("Y13").Value = 5.9
$ws.Range("Z13").Value = 0.537
$ws.Range("AC13").Value = 18.6
$ws.Range("AD13").Value = 4.6
$ws.Range("AH13").Value = 3.1
$ws.Range("AI13").Value = 12.5
$ws.Range("AJ13").Value = '1967-68'
$ws.Range("AK13").Value = 'Yes'
$ws.Range("AM13").Value = 1968

